$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new data record at row 18 ---------------------------------
# (shifts the existing rows 18..61 down to 19..62)
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44379
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112035
$ws.Range("G18").Value = "Bruselas (repollito)"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 22000
$ws.Range("L18").Value = 22000
$ws.Range("M18").Value = 22000
$ws.Range("N18").Value = "`$/malla 10 kilos"
$ws.Range("O18").Value = "Provincia de Quillota"
$ws.Range("P18").Value = 2200
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = "Hortaliza"

# --- Insert new data record at row 38 (post first insert numbering) ---
# (shifts the existing rows 38..62 down to 39..63)
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44365
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112035
$ws.Range("G38").Value = "Bruselas (repollito)"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 85
$ws.Range("K38").Value = 22000
$ws.Range("L38").Value = 22000
$ws.Range("M38").Value = 22000
$ws.Range("N38").Value = "`$/malla 10 kilos"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 2200
$ws.Range("Q38").Value = 10
$ws.Range("R38").Value = "Hortaliza"
